$d = $word.ActiveDocument

# The first line currently reads "Lab1- express app". Replace the word
# "express" with "Lab1-" (file location placeholder), turning it into
# "Lab1- Lab1- app".
$rng = $d.Content
$found = $rng.Find.Execute("express")

if ($found) {
    $rng.Text = "Lab1-"

    # Word automatically re-anchors the hidden "_GoBack" bookmark at the
    # location of the most recent edit. Collapse the range to its end
    # (right after the newly inserted text) and move the bookmark there -
    # this removes it from its previous location (in the "Lab 4" line)
    # and places it after "Lab1- Lab1-" on the first line.
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng) | Out-Null
}
